$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MYO")

# Insert two new columns before column D. Existing D:K data (and formatting)
# shifts right to F:M automatically.
$ws.Range("D1:E1").EntireColumn.Insert()

# ---- New quarter data for the two newly inserted columns (D = 31-Dec-2018, E = 30-Sep-2018) ----

# Row 7: Period Ending header dates
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D7:E7").NumberFormat = $ws.Range("F7").NumberFormat

$data = @{
  8  = @(900, 600);
  9  = @(200, 200);
  10 = @(700, 400);
  12 = @(500, 400);
  13 = @(0, 0);
  14 = @(0, "NA");
  15 = @(0, 0);
  17 = @(3600, 3300);
  18 = @(-2700, -2700);
  20 = @(0, 100);
  21 = @(-2700, -2600);
  22 = @(0, 0);
  23 = @(-2700, -2600);
  24 = @(0, 0);
  25 = @(0, 0);
  26 = @(-2700, -2600);
  27 = @(-2700, -2600);
  28 = @(0, 0);
  29 = @(0, 0);
  30 = @(0, 0);
  31 = @(0, 0);
  32 = @(0, -100);
  33 = @(-2700, -2600);
  34 = @(0, 0);
  35 = @(-2700, -2600);
}

foreach ($r in $data.Keys) {
  $vals = $data[$r]
  $dCell = $ws.Range("D" + $r)
  $eCell = $ws.Range("E" + $r)
  $dCell.NumberFormat = $ws.Range("F" + $r).NumberFormat
  $eCell.NumberFormat = $ws.Range("F" + $r).NumberFormat
  if ($vals[0] -eq "NA") { $dCell.Value2 = "NA" } else { $dCell.Value2 = $vals[0] }
  if ($vals[1] -eq "NA") { $eCell.Value2 = "NA" } else { $eCell.Value2 = $vals[1] }
}

# Row 38: Period Ending header dates (Balance Sheet)
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D38:E38").NumberFormat = $ws.Range("F38").NumberFormat

$data2 = @{
  41 = @(6500, 9100);
  42 = @(0, 0);
  43 = @(400, 400);
  44 = @(300, 300);
  45 = @(700, 700);
  46 = @(7900, 10400);
  47 = @(0, 0);
  48 = @(200, 200);
  49 = @(0, 0);
  50 = @(0, 0);
  51 = @(0, 0);
  52 = @(200, 200);
  53 = @(0, 0);
  54 = @(8300, 10800);
  57 = @(400, 400);
  58 = @("NA", "NA");
  59 = @(1400, 1400);
  60 = @(1900, 1900);
  61 = @(0, 0);
  62 = @(0, 0);
  63 = @(0, 0);
  64 = @(0, 0);
  65 = @(0, 0);
  66 = @(1900, 1900);
  68 = @(0, 0);
  69 = @(0, 0);
  70 = @(0, 0);
  71 = @(0, 0);
  72 = @(-45300, -42600);
  73 = @(0, 0);
  74 = @(0, 0);
  75 = @(0, 0);
  76 = @(6400, 8900);
  77 = @(0, 0);
}

foreach ($r in $data2.Keys) {
  $vals = $data2[$r]
  $dCell = $ws.Range("D" + $r)
  $eCell = $ws.Range("E" + $r)
  $dCell.NumberFormat = $ws.Range("F" + $r).NumberFormat
  $eCell.NumberFormat = $ws.Range("F" + $r).NumberFormat
  if ($vals[0] -eq "NA") { $dCell.Value2 = "NA" } else { $dCell.Value2 = $vals[0] }
  if ($vals[1] -eq "NA") { $eCell.Value2 = "NA" } else { $eCell.Value2 = $vals[1] }
}

# Row 80: Period Ending header dates (Cash Flow)
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D80:E80").NumberFormat = $ws.Range("F80").NumberFormat

$data3 = @{
  81  = @(-2700, -2600);
  83  = @(0, 0);
  84  = @(0, 0);
  85  = @(0, 0);
  86  = @(0, 0);
  87  = @(0, 0);
  88  = @(0, 0);
  89  = @(-2500, -2500);
  91  = @(0, 0);
  92  = @(0, 0);
  93  = @(0, 0);
  94  = @(0, 0);
  96  = @(0, 0);
  97  = @(0, 0);
  98  = @(0, 0);
  99  = @(0, 0);
  100 = @(0, -100);
  101 = @(0, 0);
  102 = @(-2600, -2600);
}

foreach ($r in $data3.Keys) {
  $vals = $data3[$r]
  $dCell = $ws.Range("D" + $r)
  $eCell = $ws.Range("E" + $r)
  $dCell.NumberFormat = $ws.Range("F" + $r).NumberFormat
  $eCell.NumberFormat = $ws.Range("F" + $r).NumberFormat
  if ($vals[0] -eq "NA") { $dCell.Value2 = "NA" } else { $dCell.Value2 = $vals[0] }
  if ($vals[1] -eq "NA") { $eCell.Value2 = "NA" } else { $eCell.Value2 = $vals[1] }
}

# Rows that are fully blank across D:M (just formatted, no values)
$blankRows = @(11, 16, 19, 39, 40, 55, 56, 67, 82, 90, 95)
foreach ($r in $blankRows) {
  $ws.Range("D" + $r + ":E" + $r).NumberFormat = $ws.Range("F" + $r).NumberFormat
}
